# feat: add 2022-Q4 data
#
# 1. "总计" (summary) sheet: the existing 2022-Q3 row becomes the 2022-Q4
#    row (new totals), and the original 2022-Q3 row is re-inserted right
#    below it with its original values.
# 2. A brand-new "2022-Q4" worksheet is inserted between "总计" and the
#    existing "2022-Q3" sheet, holding the per-fund breakdown for the new
#    quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update the "总计" summary sheet
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)

# Row 2 used to describe 2022-Q3 with 1.2; it now describes the new
# 2022-Q4 totals.
$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("D2").Value = 1.66

# Re-insert the original 2022-Q3 totals as the new row 3, matching the
# formatting already used for the index column (A2).
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("A3").Font.Bold = $true
$wsTotal.Range("A3").Borders.LineStyle = 1
$wsTotal.Range("A3").HorizontalAlignment = -4108
$wsTotal.Range("A3").VerticalAlignment = -4160

$wsTotal.Range("B3").Value = "2022-Q3"
$wsTotal.Range("C3").Value = 7
$wsTotal.Range("D3").Value = 1.2

# ---------------------------------------------------------------------
# 2. Insert a new "2022-Q4" sheet before the existing "2022-Q3" sheet
# ---------------------------------------------------------------------
$wsQ3 = $wb.Worksheets.Item(2)
$wsQ4 = $wb.Worksheets.Add($wsQ3)
$wsQ4.Name = "2022-Q4"

# Header row (bold, bordered, centered) matching the style already used
# on the other per-quarter sheets.
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $wsQ4.Cells.Item(1, $i + 2)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

# Per-fund data. Columns D/E/F/G are kept as text (e.g. "22.53") rather
# than numbers, matching the source formatting; NumberFormat = "@" forces
# that before the value is assigned.
$rows = @(
    @(0, "002446", "广发利鑫灵活配置混合A", "22.53", "73.90", "2.33", "0.5249", 9),
    @(1, "161914", "万家创业板2年定期开放混合A", "8.14", "99.81", "5.73", "0.4664", 10),
    @(2, "630008", "华商策略精选混合", "7.47", "76.58", "3.96", "0.2958", 4),
    @(3, "011172", "广发利鑫灵活配置混合C", "7.03", "73.90", "2.33", "0.1638", 9),
    @(4, "161915", "万家创业板2年定期开放混合C", "1.59", "99.81", "5.73", "0.0911", 10),
    @(5, "020015", "国泰区位优势混合A", "2.05", "84.45", "4.22", "0.0865", 8),
    @(6, "015594", "国泰区位优势混合C", "0.75", "84.45", "4.22", "0.0316", 8)
)

$r = 2
foreach ($row in $rows) {
    $a = $wsQ4.Cells.Item($r, 1)
    $a.Value = $row[0]
    $a.Font.Bold = $true
    $a.Borders.LineStyle = 1
    $a.HorizontalAlignment = -4108
    $a.VerticalAlignment = -4160

    $b = $wsQ4.Cells.Item($r, 2)
    $b.NumberFormat = "@"
    $b.Value = $row[1]

    $wsQ4.Cells.Item($r, 3).Value = $row[2]

    $d = $wsQ4.Cells.Item($r, 4)
    $d.NumberFormat = "@"
    $d.Value = $row[3]

    $e = $wsQ4.Cells.Item($r, 5)
    $e.NumberFormat = "@"
    $e.Value = $row[4]

    $f = $wsQ4.Cells.Item($r, 6)
    $f.NumberFormat = "@"
    $f.Value = $row[5]

    $g = $wsQ4.Cells.Item($r, 7)
    $g.NumberFormat = "@"
    $g.Value = $row[6]

    $wsQ4.Cells.Item($r, 8).Value = $row[7]

    $r++
}
